# Applies the VAR/arima bugfix: refreshed fitted & predicted y_values
# across the four SGP Manufacturing sheets, with one row appended to
# "y_fitted_on_begin_2016" and one row removed from "y_fitted_on_begin_2021".
$wb = $excel.ActiveWorkbook

# --- Sheet 1: y_fitted_on_begin_2016 ---
$ws1 = $wb.Worksheets.Item(1)
$data1 = @(
    @(1981, 25.82611667163626),
    @(1982, 25.25376356233554),
    @(1983, 22.48494089644078),
    @(1984, 21.85379025236798),
    @(1985, 21.80596197798535),
    @(1986, 20.8644302343111),
    @(1987, 22.94402317990297),
    @(1988, 24.19376341154965),
    @(1989, 25.82119229750707),
    @(1990, 25.22063234367931),
    @(1991, 24.44887734634226),
    @(1992, 25.22339358965022),
    @(1993, 24.26063059342838),
    @(1994, 24.28016505866036),
    @(1995, 23.50746376191397),
    @(1996, 23.69613377071189),
    @(1997, 22.8686002367677),
    @(1998, 21.94679527486855),
    @(1999, 22.28494486823695),
    @(2000, 22.46974671362737),
    @(2001, 24.86185613458518),
    @(2002, 22.77180358860602),
    @(2003, 23.82210090175429),
    @(2004, 24.02617461274749),
    @(2005, 25.90244077413681),
    @(2006, 25.83544799578697),
    @(2007, 25.42125910506903),
    @(2008, 23.16329310033411),
    @(2009, 20.5155387532922),
    @(2010, 20.26813031319024),
    @(2011, 20.71767949987342),
    @(2012, 19.63535923959325),
    @(2013, 19.09955361745289),
    @(2014, 17.82220166145943),
    @(2015, 18.04442091357657),
    @(2016, 18.030322586481)
)
for ($i = 0; $i -lt $data1.Length; $i++) {
    $r = $i + 2
    $ws1.Cells.Item($r, 1).Value = $data1[$i][0]
    $ws1.Cells.Item($r, 2).Value = $data1[$i][1]
}
# New row 37 needs C/D (lower/upper bound) cells too; clone the blank
# inline-string cells from row 2 so the pattern matches the rest of the sheet.
$ws1.Range("C2:D2").Copy($ws1.Range("C37"))

# --- Sheet 2: y_pred_on_2017_2021 (years unchanged, only y_value refreshed) ---
$ws2 = $wb.Worksheets.Item(2)
$bVals2 = @(17.40189610656805, 17.27359431982252, 17.1307094014572, 16.9789670085995, 16.82281262832847)
for ($i = 0; $i -lt $bVals2.Length; $i++) {
    $ws2.Cells.Item($i + 2, 2).Value = $bVals2[$i]
}

# --- Sheet 3: y_fitted_on_begin_2021 ---
$ws3 = $wb.Worksheets.Item(3)
# Drop the trailing 2021 row (43) now that the refreshed series is one row shorter.
$ws3.Rows.Item(43).Delete()
$data3 = @(
    @(1981, 26.14587489262068),
    @(1982, 25.47442206535145),
    @(1983, 22.72005952844342),
    @(1984, 22.04725035767255),
    @(1985, 21.97280251737266),
    @(1986, 20.93961799444227),
    @(1987, 22.85755008220352),
    @(1988, 23.96213586173698),
    @(1989, 25.46107553588415),
    @(1990, 24.80664042104852),
    @(1991, 24.06879049094408),
    @(1992, 24.86985156281298),
    @(1993, 23.98338426196701),
    @(1994, 24.06348014658787),
    @(1995, 23.38362021074435),
    @(1996, 23.64492206168755),
    @(1997, 22.99403667712899),
    @(1998, 22.17793027118497),
    @(1999, 22.55993174187589),
    @(2000, 22.65275979774851),
    @(2001, 24.96640634086572),
    @(2002, 23.03134617851963),
    @(2003, 24.08094106731323),
    @(2004, 24.10704311834144),
    @(2005, 25.92157098986607),
    @(2006, 25.85127539552811),
    @(2007, 25.41769938452232),
    @(2008, 23.204544845332),
    @(2009, 20.60107395135239),
    @(2010, 20.28616501788294),
    @(2011, 20.62305020921841),
    @(2012, 19.69584136008205),
    @(2013, 19.33457647633766),
    @(2014, 18.1786913549566),
    @(2015, 18.45635623651703),
    @(2016, 18.52420568509931),
    @(2017, 18.09331344643934),
    @(2018, 18.99110414401529),
    @(2019, 20.74471378941232),
    @(2020, 19.79821939410766),
    @(2021, 20.28030986344507)
)
for ($i = 0; $i -lt $data3.Length; $i++) {
    $r = $i + 2
    $ws3.Cells.Item($r, 1).Value = $data3[$i][0]
    $ws3.Cells.Item($r, 2).Value = $data3[$i][1]
}

# --- Sheet 4: y_pred_on_2022_2026 (years unchanged, only y_value refreshed) ---
$ws4 = $wb.Worksheets.Item(4)
$bVals4 = @(21.23613201160742, 21.45074521373596, 21.71922167024493, 22.04321072512655, 22.42469571803534)
for ($i = 0; $i -lt $bVals4.Length; $i++) {
    $ws4.Cells.Item($i + 2, 2).Value = $bVals4[$i]
}
